$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 242, pushing existing rows 242-305 down to 243-306
$ws.Rows.Item(242).Insert()

# Populate the newly inserted row 242 with the new record's data
$ws.Cells.Item(242, 1).Value = 7
$ws.Cells.Item(242, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(242, 3).Value = "Ñuble"
$ws.Cells.Item(242, 4).Value = 44798
$ws.Cells.Item(242, 5).Value = 16
$ws.Cells.Item(242, 6).Value = 100112008
$ws.Cells.Item(242, 7).Value = "Coliflor"
$ws.Cells.Item(242, 8).Value = "Sin especificar"
$ws.Cells.Item(242, 9).Value = "Segunda"
$ws.Cells.Item(242, 10).Value = 150
$ws.Cells.Item(242, 11).Value = 800
$ws.Cells.Item(242, 12).Value = 800
$ws.Cells.Item(242, 13).Value = 800
$ws.Cells.Item(242, 14).Value = "$/unidad"
$ws.Cells.Item(242, 15).Value = "Región del Maule"
$ws.Cells.Item(242, 16).Value = 800
$ws.Cells.Item(242, 17).Value = 1
$ws.Cells.Item(242, 18).Value = "Hortaliza"
